$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name
$ws.Range("C3").Value = "Parneet kaur"

# Column E (Condition being Tested) top to bottom
$ws.Range("E7").Value = "No account exists. Create with valid inputs."
$ws.Range("E8").Value = "Account created with invalid fee input"
$ws.Range("E9").Value = "Account created with date_created=2010-01-01, management_fee=2.00"
$ws.Range("E10").Value = "Account created with date_created = (today - 10 years), management_fee=2.00"
$ws.Range("E11").Value = "Account created with date_created=2023-01-01, management_fee=2.00"
$ws.Range("E12").Value = "Account created with date_created=2010-01-01, management_fee=2.00"
$ws.Range("E13").Value = "Account created with date_created=2023-01-01, management_fee=2.00"

# Column F (Method Inputs) bottom to top for rows 12-13 first, then 7-9
$ws.Range("F13").Value = "Call str(account)"
$ws.Range("F12").Value = "Call str(account"
$ws.Range("F9").Value = "Call get_service_charges()"
$ws.Range("F10").Value = "Call get_service_charges()"
$ws.Range("F11").Value = "Call get_service_charges()"
$ws.Range("F8").Value = 'management_fee="invalid"'
$ws.Range("F7").Value = "account_number=1111, client_number=22, balance=5000.00, date_created=2020-01-01, management_fee=2.00"

# Column G (Expected Result)
$ws.Range("G7").Value = "All attributes set correctly. management_fee=2.00"
$ws.Range("G8").Value = "management_fee defaults to 2.55"
$ws.Range("G9").Value = "Expected = 0.50"
$ws.Range("G10").Value = "Expected = 0.50"
$ws.Range("G11").Value = "Expected = 2.50"
$ws.Range("G12").Value = "Output includes “Management Fee: Waived”"
$ws.Range("G13").Value = "Output includes “Management Fee: $2.00”"

# Final selection, matching the saved cursor position in the workbook
[void]$ws.Range("G13").Select()
